$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for the "380 kV" case (rows 2-25, columns B,C,E,F,G,H,J,M,N,O)
$data = @{
    2 = @{ "B" = 10.39117461731422; "C" = 8.317883524737784; "E" = 23.03951531095922; "F" = 36.82758397468936; "G" = 20.34419292159597; "H" = 12.08843643360422; "J" = 7.348877636070919; "M" = 18.56271337980564; "N" = 16.61247596062897; "O" = 17.27893486685684 }
    3 = @{ "B" = 9.857326716910977; "C" = 8.069778726112093; "E" = 23.00527293212909; "F" = 36.75707984753107; "G" = 20.3145332716776; "H" = 12.12981670480522; "J" = 7.361670389933702; "M" = 18.37028534395292; "N" = 16.65169587848784; "O" = 17.33265481998719 }
    4 = @{ "B" = 9.514681102926147; "C" = 7.912623121321638; "E" = 22.98880943001238; "F" = 36.72372745437881; "G" = 20.30536942283265; "H" = 12.15744190750844; "J" = 7.369937536379401; "M" = 18.2540061683681; "N" = 16.67751064619308; "O" = 17.37023428841186 }
    5 = @{ "B" = 9.37146109769682; "C" = 7.847438543231799; "E" = 22.98325204134785; "F" = 36.7126440172935; "G" = 20.30390907752288; "H" = 12.16925648022421; "J" = 7.373410441096095; "M" = 18.20713869933065; "N" = 16.68846697308955; "O" = 17.38669979194559 }
    6 = @{ "B" = 9.347467591626762; "C" = 7.836547834657918; "E" = 22.98239891409746; "F" = 36.71095530949572; "G" = 20.30380383998577; "H" = 12.17125191118544; "J" = 7.373993402617702; "M" = 18.19938896943066; "N" = 16.69031265247306; "O" = 17.38950330869871 }
    7 = @{ "B" = 9.512763905488676; "C" = 7.911748549021838; "E" = 22.98872981282454; "F" = 36.72356781467022; "G" = 20.30534052450735; "H" = 12.15759898810872; "J" = 7.369983951812145; "M" = 18.25337194241201; "N" = 16.67765663837242; "O" = 17.37045169068202 }
    8 = @{ "B" = 10.21026898088613; "C" = 8.233376698041285; "E" = 23.02676482170883; "F" = 36.80121785152778; "G" = 20.33208886564789; "H" = 12.10224377636068; "J" = 7.353203169635904; "M" = 18.49600367629706; "N" = 16.62563958816846; "O" = 17.29650147095884 }
    9 = @{ "B" = 11.45489440535548; "C" = 8.823117961658031; "E" = 23.13730702536034; "F" = 37.03182542556395; "G" = 20.45622133966541; "H" = 12.01131158722143; "J" = 7.323554630784283; "M" = 18.98458590326235; "N" = 16.5373596749276; "O" = 17.18811312181383 }
    10 = @{ "B" = 12.28861303022919; "C" = 9.228161830749499; "E" = 23.24007291847895; "F" = 37.24815758319978; "G" = 20.59075362299903; "H" = 11.95527716652763; "J" = 7.303739675594058; "M" = 19.34852984189202; "N" = 16.48082756182047; "O" = 17.13102398777377 }
    11 = @{ "B" = 12.6495209573104; "C" = 9.405714891441564; "E" = 23.29140698972332; "F" = 37.35652370362233; "G" = 20.66121821832758; "H" = 11.93213179037835; "J" = 7.295148661649129; "M" = 19.5145802621172; "N" = 16.45690954520266; "O" = 17.10998845312477 }
    12 = @{ "B" = 12.78349465598071; "C" = 9.471944008987828; "E" = 23.31149566651563; "F" = 37.39896707578099; "G" = 20.68921637034959; "H" = 11.92370501057988; "J" = 7.29195598469541; "M" = 19.57747995292796; "N" = 16.44811046099544; "O" = 17.10273529938218 }
    13 = @{ "B" = 12.75476161213737; "C" = 9.457725767999472; "E" = 23.3071404849642; "F" = 37.389763947335; "G" = 20.68312830681003; "H" = 11.9255048309544; "J" = 7.292640895784356; "M" = 19.56393328669864; "N" = 16.44999402782719; "O" = 17.10426566133474 }
    14 = @{ "B" = 12.66059732786187; "C" = 9.411183996350911; "E" = 23.29304673773034; "F" = 37.35998748400603; "G" = 20.66349540552524; "H" = 11.93143173890131; "J" = 7.294884786093046; "M" = 19.51975494367354; "N" = 16.45618046838733; "O" = 17.10937743547103 }
    15 = @{ "B" = 12.60256666419202; "C" = 9.382543579025764; "E" = 23.28449820027853; "F" = 37.34193107917317; "G" = 20.65164032968962; "H" = 11.9351061616329; "J" = 7.296267111878597; "M" = 19.49269554953616; "N" = 16.4600034458215; "O" = 17.11260142181087 }
    16 = @{ "B" = 12.26465255521536; "C" = 9.21641969520841; "E" = 23.23680956718006; "F" = 37.24127392542225; "G" = 20.58633360331963; "H" = 11.9568370154758; "J" = 7.304309606270946; "M" = 19.33768359197113; "N" = 16.48242681191239; "O" = 17.13249826603492 }
    17 = @{ "B" = 12.05260942362938; "C" = 9.112759457704469; "E" = 23.2087218516667; "F" = 37.18205746763905; "G" = 20.54863156233036; "H" = 11.9707692724787; "J" = 7.309351547299116; "M" = 19.24267825364165; "N" = 16.49664315351423; "O" = 17.14597040639848 }
    18 = @{ "B" = 11.92892268153136; "C" = 9.052508095272877; "E" = 23.19299877463904; "F" = 37.14893638673806; "G" = 20.5278196834686; "H" = 11.97900341809392; "J" = 7.312291359987785; "M" = 19.18808285403203; "N" = 16.50498935095391; "O" = 17.15418353522391 }
    19 = @{ "B" = 11.88675003765231; "C" = 9.03200140419343; "E" = 23.18774973228472; "F" = 37.13788406958968; "G" = 20.52092362007176; "H" = 11.98182923968012; "J" = 7.313293577408437; "M" = 19.16960778290789; "N" = 16.5078443288681; "O" = 17.15704401008216 }
    20 = @{ "B" = 12.07536075514237; "C" = 9.123859661069769; "E" = 23.21166717061212; "F" = 37.18826418012452; "G" = 20.55255473922176; "H" = 11.96926331549408; "J" = 7.308810704577634; "M" = 19.25278701452862; "N" = 16.49511227797047; "O" = 17.1444881972964 }
    21 = @{ "B" = 12.68832916704766; "C" = 9.424882070399134; "E" = 23.2971688656507; "F" = 37.36869555847857; "G" = 20.66922653836478; "H" = 11.9296816873669; "J" = 7.294224059772408; "M" = 19.53273105038671; "N" = 16.45435635867893; "O" = 17.10785662315694 }
    22 = @{ "B" = 13.08214000320447; "C" = 9.615735267495936; "E" = 23.35683029086144; "F" = 37.49480891196389; "G" = 20.75313129937885; "H" = 11.90578248005611; "J" = 7.285043677609996; "M" = 19.71578433685524; "N" = 16.42922447877902; "O" = 17.08807018885352 }
    23 = @{ "B" = 12.86924907847837; "C" = 9.514424378220607; "E" = 23.32464533685711; "F" = 37.42675893198719; "G" = 20.70765610111025; "H" = 11.91835749443359; "J" = 7.289911221443717; "M" = 19.618093285571; "N" = 16.4425003400836; "O" = 17.09824955578968 }
    24 = @{ "B" = 12.06508041919202; "C" = 9.118843298282018; "E" = 23.21033426755593; "F" = 37.18545524772323; "G" = 20.55077837982305; "H" = 11.96994345998955; "J" = 7.309055091422244; "M" = 19.24821676373601; "N" = 16.49580384777003; "O" = 17.14515684703005 }
    25 = @{ "B" = 11.13202719995344; "C" = 8.668322021159954; "E" = 23.10358776896395; "F" = 36.96113127684684; "G" = 20.41499118756521; "H" = 12.03402151825649; "J" = 7.331228435469709; "M" = 18.85133231678472; "N" = 16.55977655009806; "O" = 17.21349061619719 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
